$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting the refreshed cryptos snapshot.
# A couple of D-column prices (e.g. '173.00', '0.0976') are digit-only
# strings that Excel would otherwise auto-convert to a Double (dropping
# trailing zeros / introducing float noise); those get a leading "'" so
# Excel keeps them as literal text, exactly like typing them in by hand.

$ws.Range('D2').Value = '54.913.44'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').Value = '2.297.17'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''508.51'
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('D6').Value = '''129.83'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('E7').Value = '  -0.40%  '
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('D9').Value = '2.322.94'
$ws.Range('E9').Value = '  +1.27%  '
$ws.Range('D10').Value = '''0.0976'
$ws.Range('E10').Value = '  +1.88%  '
$ws.Range('E11').Value = '  +1.71%  '
$ws.Range('D12').Value = '''5.05'
$ws.Range('E12').Value = '  +6.55%  '
$ws.Range('D13').Value = '''0.341'
$ws.Range('E13').Value = '  +1.89%  '
$ws.Range('D14').Value = '24.02'
$ws.Range('E14').Value = '  +4.71%  '
$ws.Range('D15').Value = '2.705.19'
$ws.Range('E15').Value = '  +0.59%  '
$ws.Range('D16').Value = '54.902.32'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('E17').Value = '  +1.60%  '
$ws.Range('D18').Value = '2.303.00'
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').Value = '10.73'
$ws.Range('E19').Value = '  +4.01%  '
$ws.Range('E20').Value = '  +1.19%  '
$ws.Range('D21').Value = '''6.69'
$ws.Range('E21').Value = '  +4.12%  '
$ws.Range('D22').Value = '''310.57'
$ws.Range('E22').Value = '  +1.59%  '
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D25').Value = '''0.991'
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('D27').Value = '7.54'
$ws.Range('E27').Value = '  +2.65%  '
$ws.Range('D28').Value = '''173.00'
$ws.Range('E28').Value = '  -0.77%  '
$ws.Range('D29').Value = '''6.16'
$ws.Range('E29').Value = '  +2.33%  '
$ws.Range('D30').Value = '0.0₃0712'
$ws.Range('E30').Value = '  +2.73%  '
$ws.Range('E31').Value = '  +5.35%  '
$ws.Range('E32').Value = '  +0.29%  '
$ws.Range('D33').Value = '''18.10'
$ws.Range('E33').Value = '  +1.61%  '
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('B36').Value = 'SuiNetwork'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D36').Value = '''0.921'
$ws.Range('E36').Value = '  -4.63%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '1.23'
$ws.Range('E37').Value = '  +2.58%  '
$ws.Range('E38').Value = '  +3.01%  '
$ws.Range('D39').Value = '''36.78'
$ws.Range('E39').Value = '  +2.03%  '
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').Value = '0.378'
$ws.Range('E40').Value = '  +1.12%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '''1.44'
$ws.Range('E41').Value = '  +2.14%  '
$ws.Range('D42').Value = '''135.70'
$ws.Range('E42').Value = '  +8.59%  '
$ws.Range('D43').Value = '''5.12'
$ws.Range('E43').Value = '  +5.31%  '
$ws.Range('E44').Value = '  +1.29%  '
$ws.Range('D45').Value = '''257.73'
$ws.Range('E45').Value = '  +5.79%  '
$ws.Range('D46').Value = '''0.0504'
$ws.Range('E46').Value = '  +1.46%  '
$ws.Range('E47').Value = '  +1.89%  '
$ws.Range('E48').Value = '  +1.09%  '
$ws.Range('E49').Value = '  +1.27%  '
$ws.Range('E50').Value = '  +1.29%  '
$ws.Range('E51').Value = '  +0.31%  '
